$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 19: refactoring entry for RandomVariableDifferentiableFactory ---
# Match formatting used by the analogous row above it (row 16), which has
# the same shape: package/interface/implementation/remark.
$ws.Range("A16:F16").Copy()
$ws.Range("A19:F19").PasteSpecial(-4122)

$ws.Range("B19").Value = "net.finmath.montecarlo"
$ws.Range("D19").Value = "RandomVariableDifferentiableFactory"
$ws.Range("C19").Value = "AbstractRandomVariableDifferentiableFactory"
$ws.Range("E19").Value = "interface"
$ws.Range("F19").Value = "(extracted interface) since 4.1.1"

# --- Header cell D1: "4.1.x" -> "4.1.1" ---
# Leading apostrophe forces a text value (quote-prefixed), matching the
# original author's intent of keeping a version string like "4.0.x" as text.
$ws.Range("D1").Value = "'4.1.1"
